$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("N6").Value = 2.07
$ws.Range("O6").Value = 1.83

# Row 7
$ws.Range("N7").Value = 1.5
$ws.Range("O7").Value = 2.63

# Row 12
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 3.05
$ws.Range("I12").Value = 2.2
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.35
$ws.Range("Q12").Value = 2.25
$ws.Range("R12").Value = 2.02
$ws.Range("T12").Value = 7.6
$ws.Range("U12").Value = 15
$ws.Range("V12").Value = 12
$ws.Range("X12").Value = 35
$ws.Range("Z12").Value = 6.8
$ws.Range("AA12").Value = 6.1
$ws.Range("AB12").Value = 18.5
$ws.Range("AE12").Value = 5.8
$ws.Range("AF12").Value = 9.25
$ws.Range("AG12").Value = 9.75
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 22
$ws.Range("AJ12").Value = 45

# Row 13
$ws.Range("L13").Value = 1.36
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2.1
$ws.Range("O13").Value = 1.7
